$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (23-nov) before column DX (col 128) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Columns.Item(128).Insert()

# New header cell for the inserted column
$wsPrix.Cells.Item(1, 128).Value = "23-nov"

# Fill the inserted column's data rows (2-25) with the "-" placeholder used
# throughout the sheet for missing data points.
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 128).Value = "-"
}

# --- Sheet "Gaz": append the new daily row ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A157").NumberFormat = "@"
$wsGaz.Range("A157").Value = "2025-11-21"
$wsGaz.Range("A157").ClearFormats()
$wsGaz.Range("B157").Value = 29.35

# --- Sheet "CO2": append the new daily row ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A157").NumberFormat = "@"
$wsCO2.Range("A157").Value = "2025-11-21"
$wsCO2.Range("A157").ClearFormats()
$wsCO2.Range("B157").Value = 80.28
